$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.921.96"
$ws.Range("E2").Value = "  +0.30%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.646.99"
$ws.Range("E3").Value = "  +0.73%  "

# Row 4
$ws.Range("E4").Value = "  +0.62%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.47"
$ws.Range("E5").Value = "  -0.09%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5084"
$ws.Range("E6").Value = "  +0.88%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.006"
$ws.Range("E7").Value = "  +0.37%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2575"
$ws.Range("E8").Value = "  +0.01%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06415"
$ws.Range("E9").Value = "  -0.08%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.69"
$ws.Range("E10").Value = "  +0.14%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07769"
$ws.Range("E11").Value = "  +0.76%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.688.39"
$ws.Range("E12").Value = "  +3.20%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.310"
$ws.Range("E13").Value = "  +1.31%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5463"
$ws.Range("E14").Value = "  +0.10%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅7890"
$ws.Range("E15").Value = "  -0.64%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.02"
$ws.Range("E16").Value = "  +2.35%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.001.84"
$ws.Range("E17").Value = "  +0.55%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.006"
$ws.Range("E18").Value = "  +0.43%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "197.19"
$ws.Range("E19").Value = "  -2.94%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.426"
$ws.Range("E20").Value = "  +2.25%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.03"
$ws.Range("E21").Value = "  +0.71%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.069"
$ws.Range("E22").Value = "  +1.47%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.008"
$ws.Range("E23").Value = "  +0.49%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.858"
$ws.Range("E24").Value = "  -4.51%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "140.83"
$ws.Range("E25").Value = "  -0.23%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1146"
$ws.Range("E26").Value = "  -0.21%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.890"
$ws.Range("E27").Value = "  +2.57%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.75"
$ws.Range("E28").Value = "  +0.35%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.240"
$ws.Range("E29").Value = "  -0.30%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05013"
$ws.Range("E30").Value = "  -0.32%  "

# Row 31
$ws.Range("E31").Value = "  -0.01%  "

# Row 32
$ws.Range("E32").Value = "  +0.45%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.545"
$ws.Range("E33").Value = "  +0.56%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.368"
$ws.Range("E34").Value = "  +0.35%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.8942"
$ws.Range("E35").Value = "  -0.14%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.597"
$ws.Range("E36").Value = "  -0.36%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5551"
$ws.Range("E37").Value = "  -1.22%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.132.04"
$ws.Range("E38").Value = "  -3.35%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.007"
$ws.Range("E40").Value = "  +0.50%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.666"
$ws.Range("E41").Value = "  -0.13%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8161"
$ws.Range("E42").Value = "  +1.08%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "99.74"
$ws.Range("E43").Value = "  +0.20%  "

# Row 44
$ws.Range("E44").Value = "  +7.38%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.785.16"
$ws.Range("E45").Value = "  +0.79%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4538"
$ws.Range("E46").Value = "  +0.57%  "

# Row 47
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.008"
$ws.Range("E47").Value = "  +0.40%  "

# Row 48
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "55.25"
$ws.Range("E48").Value = "  +0.75%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05093"
$ws.Range("E49").Value = "  +1.00%  "

# Row 50
$ws.Range("E50").Value = "  +0.63%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.09540"
$ws.Range("E51").Value = "  +2.62%  "
